# Weekly update: insert a new data row for "Ciboulette" (Vega Central
# Mapocho de Santiago) at row 390, pushing the existing rows 390-427
# down to 391-428. The new row carries the latest week's figures while
# every other row keeps its previous content (just shifted by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 390; this shifts rows 390:427
# down to 391:428 and copies the formatting (incl. the date number
# format in column D) from the row that used to be there.
$ws.Cells.Item(390, 1).EntireRow.Insert()

# Populate the newly inserted row 390 with the new week's data.
$ws.Cells.Item(390, 1).Value  = 9
$ws.Cells.Item(390, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(390, 3).Value  = "Metropolitana"
$ws.Cells.Item(390, 4).Value  = 44769
$ws.Cells.Item(390, 5).Value  = 13
$ws.Cells.Item(390, 6).Value  = 100112039
$ws.Cells.Item(390, 7).Value  = "Ciboulette"
$ws.Cells.Item(390, 8).Value  = "Sin especificar"
$ws.Cells.Item(390, 9).Value  = "Primera"
$ws.Cells.Item(390, 10).Value = 250
$ws.Cells.Item(390, 11).Value = 1800
$ws.Cells.Item(390, 12).Value = 2000
$ws.Cells.Item(390, 13).Value = 1900
$ws.Cells.Item(390, 14).Value = "$/docena de atados"
$ws.Cells.Item(390, 15).Value = "Región Metropolitana"
$ws.Cells.Item(390, 16).Value = 633
$ws.Cells.Item(390, 17).Value = 3
$ws.Cells.Item(390, 18).Value = "Hortaliza"
